$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("W1").Value = "Usaria o aplicativo?"
$ws.Range("W2").Value = "Não"
$ws.Range("W13").Value = "Não"
$ws.Range("W16").Value = "Não"
$ws.Range("W21").Value = "Não"
$ws.Range("W26").Value = "Sim"
$ws.Range("W31").Value = "Não"
$ws.Range("W32").Value = "Não"
$ws.Range("W33").Value = "Não"
$ws.Range("W39").Value = "Sim"
$ws.Range("W45").Value = "Sim"
$ws.Range("W46").Value = "Não"
$ws.Range("W50").Value = "Sim"
$ws.Range("W53").Value = "Não"
$ws.Range("W57").Value = "Não"
$ws.Range("W61").Value = "Sim"
$ws.Range("W64").Value = "Sim"
$ws.Range("W68").Value = "Sim"
$ws.Range("W70").Value = "Sim"
$ws.Range("W76").Value = "Não"
$ws.Range("W79").Value = "Sim"
$ws.Range("W82").Value = "Sim"
$ws.Range("W86").Value = "Sim"
$ws.Range("W93").Value = "Não"
$ws.Range("W94").Value = "Não"
$ws.Range("W99").Value = "Sim"
$ws.Range("W100").Value = "Não"
$ws.Range("W102").Value = "Sim"
$ws.Range("W108").Value = "Não"
$ws.Range("W110").Value = "Não"
$ws.Range("W115").Value = "Sim"
$ws.Range("W116").Value = "Sim"
$ws.Range("W118").Value = "Não"
$ws.Range("W123").Value = "Não"
$ws.Range("W129").Value = "Não"
$ws.Range("W134").Value = "Não"
$ws.Range("W136").Value = "Sim"
$ws.Range("W141").Value = "Não"
$ws.Range("W142").Value = "Sim"
$ws.Range("W143").Value = "Sim"
$ws.Range("W144").Value = "Não"
$ws.Range("W151").Value = "Não"
$ws.Range("W152").Value = "Não"
$ws.Range("W153").Value = "Sim"
$ws.Range("W155").Value = "Sim"
$ws.Range("W156").Value = "Sim"
$ws.Range("W158").Value = "Sim"
$ws.Range("W159").Value = "Não"
$ws.Range("W164").Value = "Sim"
$ws.Range("W166").Value = "Sim"
$ws.Range("W168").Value = "Não"
$ws.Range("W170").Value = "Não"
$ws.Range("W171").Value = "Sim"
$ws.Range("W173").Value = "Não"
$ws.Range("W179").Value = "Sim"
$ws.Range("W180").Value = "Sim"
$ws.Range("W182").Value = "Sim"
$ws.Range("W183").Value = "Sim"
$ws.Range("W184").Value = "Sim"
$ws.Range("W185").Value = "Sim"
$ws.Range("W186").Value = "Não"
$ws.Range("W189").Value = "Não"
$ws.Range("W190").Value = "Sim"
$ws.Range("W196").Value = "Sim"
$ws.Range("W197").Value = "Sim"
$ws.Range("W198").Value = "Não"
$ws.Range("W201").Value = "Sim"
$ws.Range("W205").Value = "Sim"
$ws.Range("W206").Value = "Sim"
$ws.Range("W207").Value = "Sim"
$ws.Range("W210").Value = "Não"
$ws.Range("W211").Value = "Sim"
$ws.Range("W213").Value = "Sim"
$ws.Range("W218").Value = "Não"
$ws.Range("W222").Value = "Sim"
$ws.Range("W226").Value = "Sim"
$ws.Range("W229").Value = "Sim"
$ws.Range("W230").Value = "Sim"
$ws.Range("W232").Value = "Não"
$ws.Range("W234").Value = "Não"
$ws.Range("W236").Value = "Não"
$ws.Range("W237").Value = "Não"
$ws.Range("W238").Value = "Sim"
$ws.Range("W239").Value = "Sim"
$ws.Range("W241").Value = "Não"
$ws.Range("W244").Value = "Não"
$ws.Range("W247").Value = "Não"
